$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 45 (this shifts old rows 45-141 down to 48-144)
$ws.Rows("45:47").Insert()

# Common constant values for this sheet's data rows
$A = 3
$B = "Femacal de La Calera"
$C = "Coquimbo"
$E = 5
$F = "Fruta"
$G = 100103
$H = "Frutos de hueso (carozo)"
$I = 100103002
$J = "Ciruela"

# Row 45: Black Amber / Especial
$ws.Cells.Item(45, 1).Value = $A
$ws.Cells.Item(45, 2).Value = $B
$ws.Cells.Item(45, 3).Value = $C
$ws.Cells.Item(45, 4).Value = 44581
$ws.Cells.Item(45, 5).Value = $E
$ws.Cells.Item(45, 6).Value = $F
$ws.Cells.Item(45, 7).Value = $G
$ws.Cells.Item(45, 8).Value = $H
$ws.Cells.Item(45, 9).Value = $I
$ws.Cells.Item(45, 10).Value = $J
$ws.Cells.Item(45, 11).Value = "Black Amber"
$ws.Cells.Item(45, 12).Value = "Especial"
$ws.Cells.Item(45, 13).Value = 50
$ws.Cells.Item(45, 14).Value = 13000
$ws.Cells.Item(45, 15).Value = 13000
$ws.Cells.Item(45, 16).Value = 13000
$ws.Cells.Item(45, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(45, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(45, 19).Value = 867
$ws.Cells.Item(45, 20).Value = 15

# Row 46: Black Amber / Primera
$ws.Cells.Item(46, 1).Value = $A
$ws.Cells.Item(46, 2).Value = $B
$ws.Cells.Item(46, 3).Value = $C
$ws.Cells.Item(46, 4).Value = 44581
$ws.Cells.Item(46, 5).Value = $E
$ws.Cells.Item(46, 6).Value = $F
$ws.Cells.Item(46, 7).Value = $G
$ws.Cells.Item(46, 8).Value = $H
$ws.Cells.Item(46, 9).Value = $I
$ws.Cells.Item(46, 10).Value = $J
$ws.Cells.Item(46, 11).Value = "Black Amber"
$ws.Cells.Item(46, 12).Value = "Primera"
$ws.Cells.Item(46, 13).Value = 55
$ws.Cells.Item(46, 14).Value = 12000
$ws.Cells.Item(46, 15).Value = 12000
$ws.Cells.Item(46, 16).Value = 12000
$ws.Cells.Item(46, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(46, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(46, 19).Value = 800
$ws.Cells.Item(46, 20).Value = 15

# Row 47: Black Amber / Segunda
$ws.Cells.Item(47, 1).Value = $A
$ws.Cells.Item(47, 2).Value = $B
$ws.Cells.Item(47, 3).Value = $C
$ws.Cells.Item(47, 4).Value = 44581
$ws.Cells.Item(47, 5).Value = $E
$ws.Cells.Item(47, 6).Value = $F
$ws.Cells.Item(47, 7).Value = $G
$ws.Cells.Item(47, 8).Value = $H
$ws.Cells.Item(47, 9).Value = $I
$ws.Cells.Item(47, 10).Value = $J
$ws.Cells.Item(47, 11).Value = "Black Amber"
$ws.Cells.Item(47, 12).Value = "Segunda"
$ws.Cells.Item(47, 13).Value = 50
$ws.Cells.Item(47, 14).Value = 10000
$ws.Cells.Item(47, 15).Value = 10000
$ws.Cells.Item(47, 16).Value = 10000
$ws.Cells.Item(47, 17).Value = "$/caja 15 kilos empedrada"
$ws.Cells.Item(47, 18).Value = "Provincia de San Felipe de Aconcagua"
$ws.Cells.Item(47, 19).Value = 667
$ws.Cells.Item(47, 20).Value = 15
